# data group - edit DONE
# Adds a new "katalonsepuluh" negative test-data row (row 8) to the
# "DataGroup-add" sheet, mirroring the existing rows: Action=Add,
# groupID=katalonsepuluh, customerName=All, principalName=All,
# category=NEGATIVE (deskripsi left blank, same as the row's source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataGroup-add")

$ws.Range("B8").Value = "Add"
$ws.Range("C8").Value = "katalonsepuluh"
$ws.Range("E8").Value = "All"
$ws.Range("F8").Value = "All"
$ws.Range("G8").Value = "NEGATIVE"

# Move the visible cursor to E6, matching the author's final selection.
$null = $ws.Range("E6").Select()

# Restore the application window size recorded in the author's commit.
$win = $excel.ActiveWindow
$win.Width = 19635
$win.Height = 7500
